$d = $word.ActiveDocument
$rng = $d.Content
$rng.Collapse(0)
$xmlBody = @'
<w:p>
      <w:r>
        <w:t>Il existe deux manières pour maintenir un équilibre hydrique :</w:t>
      </w:r>
    </w:p>
    <w:tbl>
      <w:tblPr>
        <w:tblStyle w:val="Grilledetableauclaire"/>
        <w:tblW w:w="0" w:type="auto"/>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="3667"/>
        <w:gridCol w:w="3668"/>
      </w:tblGrid>
      <w:tr>
        <w:trPr>
          <w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            <w:tcW w:w="3667" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Osmotolérance</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3668" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            </w:pPr>
            <w:r>
              <w:t>osmorégulation</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            <w:tcW w:w="3667" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Ayant la même osmolarité que leur environnement</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3668" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            </w:pPr>
            <w:r>
              <w:t>Régule leur osmolarité</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            <w:tcW w:w="3667" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:r>
              <w:t>Plutôt les animaux qui vivent dans un milieu aquatique stable</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3668" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
            </w:pPr>
            <w:r>
              <w:t>Concerne les animaux d’eau douce et terrestre</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Les invertébrés marins sont souvent de type osmotolérant. Seul qui </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Pour les Vertébrés marins osmorégulateur, l’océan est déshydratant.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Chez les chondrichtyen, le milieu intracellulaire est hyper-osmotique à cause de l’accumulation des déchets métaboliques comme l’urée qui ne sont pas éliminé de leur organisme. Il en résulte une entré d’eau qui est évacuée par leur organes.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Dulcicoles milieu hypotonique se qui provoque une entrer </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Et une perte de soluté qui est compensé par le gain fait par la nourriture. Possède des cellules spécialisées transport actif pour maintenir une concentration.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Saumon est capable de modifier la taille des cellules sécrétrice de sel grâce à une hormone.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Milieu aquatique précaire animaux survivre à un asséchement grâce à un état</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Hydrobiose</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> état d’inactivité ou l’organisme perd l’eau.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Animaux terrestre sont menacés par la déshydratation </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Un homme meurt lorsqu’il a perdu 12% de son eau.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Adaptation physiologique cuticule chez les Insectes.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Comportementales période d’activité la nuit </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Les pertes dû à l’évaporation, dans les selles et l’urine sont comblées par la respiration cellulaire et l’alimentation (notamment par </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Le maintien d’un milieu interne à un coût énergétique d’autant plus important qu’il est différent du milieu externe. Il dépend également :</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>De la surface d’échanges</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>De la perméabilité des membranes</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Chez les poissons, 5% de l’énergie totale.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Rmq : Au sein de groupes d’organismes l’évolution à sélectionner les individus ayant la concentration qui minimisent la différence avec le milieu ce qui contribuent </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>a</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> réduire la dépense énergétique </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Le maintien des concentrations de solutés se fait au sein de structures spécialisés. Comme les reins chez les vertébrés </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Insecte système ouvert le maintien se fait par l’hémolymphe </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Réguler diminuer en éliminant.</w:t>
      </w:r>
    </w:p>

'@

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' + $xmlBody + '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$rng.InsertXML($xml)
Write-Output "Inserted content successfully"
